$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price/date data (columns D, J, K, L, M, P) is being re-shuffled
# across rows 2-12; all other columns (A,B,C,E,F,G,H,I,N,O,Q,R) stay the same
# for every row. Write the new target values for each row.

$values = @{
    2  = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    3  = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
    4  = @{ D = 45132; J = 170;  K = 2200; L = 2500; M = 2359; P = 2359 }
    5  = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    6  = @{ D = 44907; J = 2300; K = 900;  L = 1000; M = 952;  P = 952  }
    7  = @{ D = 44895; J = 200;  K = 1200; L = 1300; M = 1255; P = 1255 }
    8  = @{ D = 45062; J = 1700; K = 2800; L = 3000; M = 2900; P = 2900 }
    9  = @{ D = 44175; J = 1400; K = 1900; L = 2000; M = 1950; P = 1950 }
    10 = @{ D = 44883; J = 290;  K = 1400; L = 1500; M = 1434; P = 1434 }
    11 = @{ D = 44893; J = 3300; K = 1200; L = 1300; M = 1261; P = 1261 }
    12 = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("J$row").Value = $rowVals.J
    $ws.Range("K$row").Value = $rowVals.K
    $ws.Range("L$row").Value = $rowVals.L
    $ws.Range("M$row").Value = $rowVals.M
    $ws.Range("P$row").Value = $rowVals.P
}
